# feat: add 2022-Q4 data
#
# The "总计" (total) summary sheet gains a new top row for the 2022-Q4
# period (pushing the existing 2021-Q2 summary row down), and a brand new
# "2022-Q4" worksheet is inserted between "总计" and "2021-Q2" holding the
# per-fund holdings detail for that quarter.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsOld   = $wb.Worksheets.Item(2)   # currently "2021-Q2"

# ---------------------------------------------------------------------
# 1) Insert a new worksheet right before the existing "2021-Q2" sheet so
#    the final tab order is 总计, 2022-Q4, 2021-Q2.
# ---------------------------------------------------------------------
$wsNew = $wb.Worksheets.Add($wsOld, $null)
$wsNew.Name = "2022-Q4"

# Carry over the header formatting (bold/border style) from the 总计
# sheet's header row so the new sheet visually matches its siblings.
$wsTotal.Range("B1:D1").Copy($wsNew.Range("B1:D1"))
$wsTotal.Range("B1:D1").Copy($wsNew.Range("E1:G1"))
$wsTotal.Range("B1").Copy($wsNew.Range("H1"))
$wsTotal.Range("A2").Copy($wsNew.Range("A2"))

# Header row text
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Data row (row 2) - fund code/name/ratios are stored as text in the
# source data, same as the sibling quarter sheets, so force a text
# number-format before assigning the numeric-looking strings.
$wsNew.Range("A2").Value = 0

$wsNew.Range("B2").NumberFormat = "@"
$wsNew.Range("B2").Value = "233009"

$wsNew.Range("C2").NumberFormat = "@"
$wsNew.Range("C2").Value = "大摩多因子精选策略混合"

$wsNew.Range("D2").NumberFormat = "@"
$wsNew.Range("D2").Value = "6.42"

$wsNew.Range("E2").NumberFormat = "@"
$wsNew.Range("E2").Value = "91.11"

$wsNew.Range("F2").NumberFormat = "@"
$wsNew.Range("F2").Value = "0.94"

$wsNew.Range("G2").NumberFormat = "@"
$wsNew.Range("G2").Value = "0.0603"

$wsNew.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: push the existing 2021-Q2 summary row from
#    row 2 down to row 3 (keeping its formatting), then write the new
#    2022-Q4 summary values into row 2.
# ---------------------------------------------------------------------
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.06
